$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'New 닥터퓨리 간편청소 변기세정제 스틱+청소포 12개입 포함 변기클리너스틱'
$ws.Range("C33").Value = 'https://naver.me/FLecmA7I'
$ws.Range("D33").Value = 'https://shop-phinf.pstatic.net/20250822_196/1755829583922m2i9q_JPEG/7599366041362571_1796467425.jpg'
$ws.Range("E33").Value = '간편 변기 세정 스틱'
$ws.Range("F33").Value = '번거로운 변기 청소, 이 스틱 하나로 간편하게 해결하세요. 청결과 상쾌함을 동시에 느껴보세요.'

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = '행운의부 소원부적 황금 소원 거북이 개업 승진 재물 합격 삼재 새해 선물 부적'
$ws.Range("C34").Value = 'https://naver.me/5FENBwGE'
$ws.Range("D34").Value = 'https://shop-phinf.pstatic.net/20250118_39/1737189262496Akyes_JPEG/71322074642706286_584694950.jpg'
$ws.Range("E34").Value = '황금 소원 거북이 부적'
$ws.Range("F34").Value = '개업·승진·재물·합격 등 행운을 바라는 당신의 소원에 힘을 더해 드려요.'

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = '[특별행사] 짱구 벽걸이 달력 2025년 캘린더 오렌지 연말 선물 캐릭터 굿즈'
$ws.Range("C35").Value = 'https://naver.me/Fc5i9Wl8'
$ws.Range("D35").Value = 'https://shop-phinf.pstatic.net/20241025_235/1729821965410is2tP_JPEG/6309273534817397_718649002.jpg'
$ws.Range("E35").Value = '2025년 짱구 벽걸이 달력'
$ws.Range("F35").Value = '귀여운 짱구와 함께하는 새해! 일정을 한눈에 관리하며 특별한 연말을 완성해보세요.'

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = '데이로이 미니 복주머니 전통 복조리 새해 복주머니 외국인 기념품 명절 선물'
$ws.Range("C36").Value = 'https://naver.me/G9pWLJG3'
$ws.Range("D36").Value = 'https://shop-phinf.pstatic.net/20251117_226/17633852195067Mw1T_JPEG/18065954483823832_2026591829.jpg'
$ws.Range("E36").Value = '작은 복주머니 선물용'
$ws.Range("F36").Value = '소중한 분께 전통의 의미를 담은 미니 복주머니로 새해 행운을 전해보세요.'

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = '메리크리스마스카드 8cm 30개입 A,B형 택1 해피뉴이어 새해 감사 성탄 미니카드'
$ws.Range("C37").Value = 'https://naver.me/FivD9xg5'
$ws.Range("D37").Value = 'https://shop-phinf.pstatic.net/20230602_147/1685682054414QwpOL_JPEG/31425722217311814_1865601803.jpg'
$ws.Range("E37").Value = '미니 크리스마스 감사카드'
$ws.Range("F37").Value = '소중한 마음을 전하는 작은 카드, 한 장으로 따뜻한 인사를 완성해 보세요.'

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = '돈쓸어담는 황금빗자루 각인제작 액막이 신혼집들이 개업 신년새해 이사선물'
$ws.Range("C38").Value = 'https://naver.me/xxFliNHt'
$ws.Range("D38").Value = 'https://shop-phinf.pstatic.net/20250924_270/1758705738235S33tL_JPEG/8796685156155698_856278412.jpg'
$ws.Range("E38").Value = '황금빗자루 액막이 선물'
$ws.Range("F38").Value = '새 출발을 축하하는 황금빗자루로 액운을 쓸어내고 행운을 담아보세요.'

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = '2024 갑진년 용의 해 신년 행운 푸른 용 용띠키링'
$ws.Range("C39").Value = 'https://naver.me/xv64J5Mp'
$ws.Range("D39").Value = 'https://shop-phinf.pstatic.net/20251207_17/1765062713833zPYpN_JPEG/745897979639105_606870278.jpg'
$ws.Range("E39").Value = '2024 신년 행운 푸른 용 키링'
$ws.Range("F39").Value = '새해 시작에 행운을 더해줄 용띠 맞춤 키링으로 특별한 기운을 느껴보세요.'

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = '2026년 달력 2026 캘린더 벽걸이달력 은행 카렌다 새해'
$ws.Range("C40").Value = 'https://naver.me/IM4Y7ArS'
$ws.Range("D40").Value = 'https://shop-phinf.pstatic.net/20251028_82/1761631215994YUn7H_JPEG/11301190108030554_319270666.jpg'
$ws.Range("E40").Value = '2026년 벽걸이 달력'
$ws.Range("F40").Value = '새해 시작을 깔끔하게! 2026년 한눈에 확인하고 계획 세우세요.'
